# Update the Categories values to use textual keys instead of numeric IDs,
# and set the (previously empty/numeric) Categories cell for the bicycle
# product row to "Bicycles".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Cars,Pick-Ups "
$ws.Range("E3").Value = "Cars,Roadsters"
$ws.Range("E4").Value = "Bicycles"

$ws.Range("F4").Select()
